$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the contents of columns L and M (header row 1 through data row 16)
for ($r = 1; $r -le 16; $r++) {
    $lCell = $ws.Cells.Item($r, 12)  # column L
    $mCell = $ws.Cells.Item($r, 13)  # column M

    $lVal = $lCell.Value2
    $mVal = $mCell.Value2

    $lCell.Value = $mVal
    $mCell.Value = $lVal
}
